$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the lab end time for the second row (D2) from "6:00PM" to "6:45PM"
$ws.Range("D2").Value = "6:45PM"

# Move the active selection to D2 (matches the edited cell)
$ws.Range("D2").Select()
